$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOV-2020")

# ---------------------------------------------------------------
# 1) Fix task-number values in rows 4 and 5 (A4: 2 -> 3, A5: 2 -> 4)
# ---------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# ---------------------------------------------------------------
# 2) Capture formatting templates BEFORE they get overwritten.
#    Row 7 (A7:G7) currently carries the "empty template" styling
#    that rows 8 and 9 need to inherit (F column style 39).
# ---------------------------------------------------------------
$ws.Range("A7:G7").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122)
$ws.Range("A9:G9").PasteSpecial(-4122)

# Row 5 (A5:G5) holds the fully "active" styling (incl. D col style 14
# and F col style 22) that rows 6, 7 and 10 need.
$ws.Range("A5:G5").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)
$ws.Range("A7:G7").PasteSpecial(-4122)
$ws.Range("A10:G10").PasteSpecial(-4122)

# D column (wrap/left/vcenter) style for rows 8 and 9 should match D5/D6 (style 14)
$ws.Range("D5").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)

# F10 ("COMPLETED") needs the green-ish "Completed" style that already
# lives on sheet "SEP-2020" (F15) in this workbook.
$wsSep = $wb.Worksheets.Item("SEP-2020")
$wsSep.Range("F15").Copy()
$ws.Range("F10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 3) Fill in the new day rows (6, 7, 8, 9) and the completed row (10)
# ---------------------------------------------------------------

# Nov 5, 2020 - nMVAR_QA, WIP, 50%
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 44140
$ws.Range("C6").Value = "nMVAR "
$ws.Range("D6").Value = "nMVAR_QA           "
$ws.Range("E6").Value = 0.5
$ws.Range("F6").Value = "WIP"

# Nov 6, 2020 - nMVAR_QA, WIP, 80%
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 44141
$ws.Range("C7").Value = "nMVAR "
$ws.Range("D7").Value = "nMVAR_QA           "
$ws.Range("E7").Value = 0.8
$ws.Range("F7").Value = "WIP"

# Nov 7, 2020 - number/date only, rest left blank
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 44142

# Nov 8, 2020 - number/date only, rest left blank
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 44143

# Nov 9, 2020 - nMVAR_QA / nMVAR_TSS, COMPLETED, 100%/10%
# (status marked COMPLETED first, then the task text, then the % split -
#  this mirrors the order the original strings were added to the workbook)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 44144
$ws.Range("C10").Value = "nMVAR "
$ws.Range("F10").Value = "COMPLETED"
$ws.Range("D10").Value = "1) nMVAR_QA                                                                     2)nMVAR_TSS       "
$ws.Range("E10").Value = "1) 100% 2)10%"

$ws.Rows.Item(10).RowHeight = 29.4

# ---------------------------------------------------------------
# 4) Column D got narrower after the new text was entered
# ---------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 37

# ---------------------------------------------------------------
# 5) Selection moved to F13 before the file was last saved
# ---------------------------------------------------------------
$ws.Range("F13").Select()
